$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 523.63635
$ws.Range("I2").Value = 169.875
$ws.Range("J2").Value = 1467
$ws.Range("K2").Value = 169.875
$ws.Range("L2").Value = 1467
$ws.Range("M2").Value = -56.875
$ws.Range("N2").Value = -1693
$ws.Range("H40").Value = 2249.818
$ws.Range("I40").Value = 2128.5715
$ws.Range("J40").Value = 2306.4
$ws.Range("K40").Value = 2128.5715
$ws.Range("L40").Value = 2306.4
$ws.Range("M40").Value = -1953.5715
$ws.Range("N40").Value = -2656.4
$ws.Range("H97").Value = 2803.6667
$ws.Range("J97").Value = 3204.4
$ws.Range("L97").Value = 9613.200000000001
$ws.Range("N97").Value = -10605.2
$ws.Range("H112").Value = 3309.577
$ws.Range("I112").Value = 924.8333
$ws.Range("J112").Value = 4025
$ws.Range("K112").Value = 2774.4999
$ws.Range("L112").Value = 12075
$ws.Range("M112").Value = -1666.4999
$ws.Range("N112").Value = -14291
$ws.Range("H135").Value = 1652.5555
$ws.Range("I135").Value = 1454.7858
$ws.Range("J135").Value = 2344.75
$ws.Range("K135").Value = 13093.0722
$ws.Range("L135").Value = 21102.75
$ws.Range("M135").Value = -10558.0722
$ws.Range("N135").Value = -26172.75
$ws.Range("H137").Value = 1804.3
$ws.Range("I137").Value = 1586.4419
$ws.Range("K137").Value = 4759.3257
$ws.Range("M137").Value = -2209.3257
$ws.Range("H139").Value = 90197.836
$ws.Range("J139").Value = 98095.60000000001
$ws.Range("L139").Value = 98095.60000000001
$ws.Range("N139").Value = -108375.6
$ws.Range("H141").Value = 4988.2666
$ws.Range("I141").Value = 1935.4445
$ws.Range("J141").Value = 9567.5
$ws.Range("K141").Value = 5806.333500000001
$ws.Range("L141").Value = 28702.5
$ws.Range("M141").Value = -626.3335000000006
$ws.Range("N141").Value = -39062.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5100.95
$ws.Range("I58").Value = 3617.1538
$ws.Range("J58").Value = 7856.5713
$ws.Range("K58").Value = 3617.1538
$ws.Range("L58").Value = 7856.5713
$ws.Range("M58").Value = -3414.1538
$ws.Range("N58").Value = -8262.5713
$ws.Range("H132").Value = 1590.5088
$ws.Range("I132").Value = 1261.96
$ws.Range("J132").Value = 3937.2856
$ws.Range("K132").Value = 3785.88
$ws.Range("L132").Value = 11811.8568
$ws.Range("M132").Value = -1255.88
$ws.Range("N132").Value = -16871.8568
$ws.Range("H136").Value = 5100.95
$ws.Range("I136").Value = 3617.1538
$ws.Range("J136").Value = 7856.5713
$ws.Range("K136").Value = 10851.4614
$ws.Range("L136").Value = 23569.7139
$ws.Range("M136").Value = -8301.4614
$ws.Range("N136").Value = -28669.7139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 335.69232
$ws.Range("J2").Value = 348.42856
$ws.Range("L2").Value = 2090.57136
$ws.Range("N2").Value = -2316.57136
$ws.Range("H19").Value = 4468
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H38").Value = 433.66666
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H69").Value = 794
$ws.Range("I69").Value = 794
$ws.Range("K69").Value = 2382
$ws.Range("M69").Value = -1571
$ws.Range("H72").Value = 794
$ws.Range("I72").Value = 794
$ws.Range("K72").Value = 7146
$ws.Range("M72").Value = -3090
$ws.Range("H116").Value = 2735.7144
$ws.Range("I116").Value = 975
$ws.Range("J116").Value = 3440
$ws.Range("K116").Value = 2925
$ws.Range("L116").Value = 10320
$ws.Range("M116").Value = 517
$ws.Range("N116").Value = -17204
$ws.Range("H131").Value = 912093.3
$ws.Range("I131").Value = 1007.0909
$ws.Range("K131").Value = 3021.2727
$ws.Range("M131").Value = 2018.7273
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9010.444
$ws.Range("I80").Value = 2728.4285
$ws.Range("J80").Value = 30997.5
$ws.Range("K80").Value = 2728.4285
$ws.Range("L80").Value = 30997.5
$ws.Range("M80").Value = -1730.4285
$ws.Range("N80").Value = -32993.5
$ws.Range("H83").Value = 9010.444
$ws.Range("I83").Value = 2728.4285
$ws.Range("J83").Value = 30997.5
$ws.Range("K83").Value = 13642.1425
$ws.Range("L83").Value = 154987.5
$ws.Range("M83").Value = -8650.1425
$ws.Range("N83").Value = -164971.5
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080
$ws.Range("H113").Value = 2186.75
$ws.Range("I113").Value = 2249.25
$ws.Range("K113").Value = 2249.25
$ws.Range("M113").Value = -79.25
$ws.Range("H114").Value = 209999.5
$ws.Range("J114").Value = 209999.5
$ws.Range("L114").Value = 209999.5
$ws.Range("N114").Value = -218677.5
$ws.Range("H126").Value = 3498.077
$ws.Range("I126").Value = 3258.3704
$ws.Range("J126").Value = 4037.4167
$ws.Range("K126").Value = 9775.111199999999
$ws.Range("L126").Value = 12112.2501
$ws.Range("M126").Value = -7305.111199999999
$ws.Range("N126").Value = -17052.2501
$ws.Range("H132").Value = 5987.684
$ws.Range("I132").Value = 6876.727
$ws.Range("K132").Value = 20630.181
$ws.Range("M132").Value = -18100.181
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6650.4614
$ws.Range("I40").Value = 4768.727
$ws.Range("J40").Value = 17000
$ws.Range("K40").Value = 4768.727
$ws.Range("L40").Value = 17000
$ws.Range("M40").Value = -4632.727
$ws.Range("N40").Value = -17272
$ws.Range("H46").Value = 3637.889
$ws.Range("I46").Value = 614.3333
$ws.Range("J46").Value = 5149.6665
$ws.Range("K46").Value = 614.3333
$ws.Range("L46").Value = 5149.6665
$ws.Range("M46").Value = -426.3333
$ws.Range("N46").Value = -5525.6665
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 12450.409
$ws.Range("I96").Value = 3481.4
$ws.Range("J96").Value = 19924.584
$ws.Range("K96").Value = 3481.4
$ws.Range("L96").Value = 19924.584
$ws.Range("M96").Value = -2108.4
$ws.Range("N96").Value = -22670.584
